$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.528.27"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "1.650.50"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").Value = "'300.52"
$ws.Range("E6").Value = "  -0.97%  "

$ws.Range("D7").Value = "'0.3790"
$ws.Range("E7").Value = "  +0.86%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3574"
$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'50.75"
$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("D10").Value = "'0.08113"
$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "'22.11"
$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").Value = "'6.416"
$ws.Range("E14").Value = "  -1.53%  "

$ws.Range("D15").Value = "'7.420"
$ws.Range("E15").Value = "  +1.21%  "

$ws.Range("D16").Value = "'0.00001204"
$ws.Range("E16").Value = "  -1.68%  "

$ws.Range("D17").Value = "1.658.49"
$ws.Range("E17").Value = "  +1.47%  "

$ws.Range("D18").Value = "'97.04"
$ws.Range("E18").Value = "  +0.83%  "

$ws.Range("D19").Value = "'0.06997"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").Value = "'6.791"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "'12.62"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("D24").Value = "23.540.54"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").Value = "'2.497"
$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("D26").Value = "'2.934"
$ws.Range("E26").Value = "  -5.78%  "

$ws.Range("D27").Value = "'20.99"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").Value = "'152.00"
$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("D29").Value = "'5.235"
$ws.Range("E29").Value = "  +0.84%  "

$ws.Range("D30").Value = "'133.04"
$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("D31").Value = "1.837.43"
$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("D32").Value = "'6.996"
$ws.Range("E32").Value = "  +5.02%  "

$ws.Range("D33").Value = "'2.148"
$ws.Range("E33").Value = "  +5.19%  "

$ws.Range("D34").Value = "'11.90"
$ws.Range("E34").Value = "  +3.69%  "

$ws.Range("D35").Value = "'1.042"
$ws.Range("E35").Value = "  -4.86%  "

$ws.Range("D36").Value = "'0.02741"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("D37").Value = "'0.08704"
$ws.Range("E37").Value = "  -0.57%  "

$ws.Range("D38").Value = "'0.2455"
$ws.Range("E38").Value = "  -1.35%  "

$ws.Range("D39").Value = "'5.994"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("D40").Value = "'13.17"
$ws.Range("E40").Value = "  +4.73%  "

$ws.Range("D41").Value = "'0.06884"
$ws.Range("E41").Value = "  -1.23%  "

$ws.Range("D42").Value = "'0.6931"
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("D44").Value = "'15.75"
$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("D45").Value = "'0.6457"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "'2.268"
$ws.Range("E47").Value = "  -2.36%  "

$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").Value = "'0.07829"
$ws.Range("E49").Value = "  -1.15%  "

$ws.Range("D50").Value = "'127.28"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").Value = "'1.176"
$ws.Range("E51").Value = "  -0.19%  "
